$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 224172
$ws.Range("J17").Value = 224172
$ws.Range("L17").Value = 672516
$ws.Range("N17").Value = -672852
$ws.Range("H38").Value = 1491.5
$ws.Range("I38").Value = 1172.5454
$ws.Range("K38").Value = 3517.6362
$ws.Range("M38").Value = -3145.6362
$ws.Range("H40").Value = 55556556
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 111111110
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 111111110
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -111111460
$ws.Range("H103").Value = 5739
$ws.Range("I103").Value = 2000
$ws.Range("J103").Value = 6673.75
$ws.Range("K103").Value = 6000
$ws.Range("L103").Value = 20021.25
$ws.Range("M103").Value = -5414
$ws.Range("N103").Value = -21193.25
$ws.Range("H107").Value = 1402.4445
$ws.Range("I107").Value = 219.33333
$ws.Range("K107").Value = 219.33333
$ws.Range("M107").Value = 1700.66667
$ws.Range("H121").Value = 1349
$ws.Range("J121").Value = 1349
$ws.Range("L121").Value = 4047
$ws.Range("N121").Value = -7541
$ws.Range("H138").Value = 3833.1956
$ws.Range("I138").Value = 2662.125
$ws.Range("J138").Value = 5110.727
$ws.Range("K138").Value = 7986.375
$ws.Range("L138").Value = 15332.181
$ws.Range("M138").Value = -2846.375
$ws.Range("N138").Value = -25612.181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 30200.2
$ws.Range("I6").Value = 17000.334
$ws.Range("K6").Value = 17000.334
$ws.Range("M6").Value = -16827.334
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H32").Value = 6701.9346
$ws.Range("I32").Value = 3527.457
$ws.Range("K32").Value = 3527.457
$ws.Range("M32").Value = -3240.457
$ws.Range("H46").Value = 3999.75
$ws.Range("I46").Value = 999.5
$ws.Range("J46").Value = 7000
$ws.Range("K46").Value = 999.5
$ws.Range("L46").Value = 7000
$ws.Range("M46").Value = -680.5
$ws.Range("N46").Value = -7638
$ws.Range("H61").Value = 252256880
$ws.Range("I61").Value = 252256880
$ws.Range("K61").Value = 252256880
$ws.Range("M61").Value = -252256668
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H136").Value = 252256880
$ws.Range("I136").Value = 252256880
$ws.Range("K136").Value = 756770640
$ws.Range("M136").Value = -756768090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2958.25
$ws.Range("I5").Value = 3313.8333
$ws.Range("J5").Value = 1891.5
$ws.Range("K5").Value = 3313.8333
$ws.Range("L5").Value = 1891.5
$ws.Range("M5").Value = -3200.8333
$ws.Range("N5").Value = -2117.5
$ws.Range("H7").Value = 641
$ws.Range("I7").Value = 664.8333
$ws.Range("K7").Value = 664.8333
$ws.Range("M7").Value = -551.8333
$ws.Range("H54").Value = 3416.3333
$ws.Range("I54").Value = 2624.5
$ws.Range("K54").Value = 2624.5
$ws.Range("M54").Value = -2140.5
$ws.Range("H99").Value = 1567.28
$ws.Range("I99").Value = 862.2143
$ws.Range("K99").Value = 862.2143
$ws.Range("M99").Value = 635.7857
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 132.09091
$ws.Range("I7").Value = 31.75
$ws.Range("K7").Value = 31.75
$ws.Range("M7").Value = 81.25
$ws.Range("H12").Value = 2766.8
$ws.Range("I12").Value = 2276.1667
$ws.Range("K12").Value = 2276.1667
$ws.Range("M12").Value = -2106.1667
$ws.Range("H16").Value = 931.5925999999999
$ws.Range("I16").Value = 757.8182
$ws.Range("K16").Value = 757.8182
$ws.Range("M16").Value = -470.8182
$ws.Range("H31").Value = 8836.125
$ws.Range("I31").Value = 2366.2104
$ws.Range("J31").Value = 18292.154
$ws.Range("K31").Value = 2366.2104
$ws.Range("L31").Value = 18292.154
$ws.Range("M31").Value = -2071.2104
$ws.Range("N31").Value = -18882.154
$ws.Range("H34").Value = 8836.125
$ws.Range("I34").Value = 2366.2104
$ws.Range("J34").Value = 18292.154
$ws.Range("K34").Value = 2366.2104
$ws.Range("L34").Value = 18292.154
$ws.Range("M34").Value = -2164.2104
$ws.Range("N34").Value = -18696.154
$ws.Range("H51").Value = 13100
$ws.Range("I51").Value = 13100
$ws.Range("K51").Value = 13100
$ws.Range("M51").Value = -12364
$ws.Range("H58").Value = 62513388
$ws.Range("I58").Value = 62513388
$ws.Range("K58").Value = 62513388
$ws.Range("M58").Value = -62513185
$ws.Range("H60").Value = 54999.5
$ws.Range("J60").Value = 54999.5
$ws.Range("L60").Value = 54999.5
$ws.Range("N60").Value = -56021.5
$ws.Range("H61").Value = 13100
$ws.Range("I61").Value = 13100
$ws.Range("K61").Value = 13100
$ws.Range("M61").Value = -12752
$ws.Range("H108").Value = 41666.668
$ws.Range("I108").Value = 20000
$ws.Range("J108").Value = 52500
$ws.Range("K108").Value = 20000
$ws.Range("L108").Value = 52500
$ws.Range("M108").Value = -16160
$ws.Range("N108").Value = -60180
$ws.Range("H113").Value = 931.5925999999999
$ws.Range("I113").Value = 757.8182
$ws.Range("K113").Value = 757.8182
$ws.Range("M113").Value = 1412.1818
$ws.Range("H136").Value = 62513388
$ws.Range("I136").Value = 62513388
$ws.Range("K136").Value = 187540164
$ws.Range("M136").Value = -187537614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 175.9
$ws.Range("J2").Value = 255.8
$ws.Range("L2").Value = 1534.8
$ws.Range("N2").Value = -1760.8
$ws.Range("H38").Value = 22.2
$ws.Range("I38").Value = 13.2
$ws.Range("J38").Value = 31.2
$ws.Range("K38").Value = 39.59999999999999
$ws.Range("L38").Value = 93.59999999999999
$ws.Range("M38").Value = 307.4
$ws.Range("N38").Value = -787.6
$ws.Range("H75").Value = 4333
$ws.Range("J75").Value = 3999.5
$ws.Range("L75").Value = 11998.5
$ws.Range("N75").Value = -13994.5
$ws.Range("H78").Value = 4333
$ws.Range("J78").Value = 3999.5
$ws.Range("L78").Value = 35995.5
$ws.Range("N78").Value = -45979.5
$ws.Range("H88").Value = 18285.572
$ws.Range("I88").Value = 9999
$ws.Range("K88").Value = 29997
$ws.Range("M88").Value = -29569
$ws.Range("H91").Value = 18285.572
$ws.Range("I91").Value = 9999
$ws.Range("K91").Value = 29997
$ws.Range("M91").Value = -28515
$ws.Range("H126").Value = 7462
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 9943
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 29829
$ws.Range("M126").Value = -2560
$ws.Range("N126").Value = -39709
$ws.Range("H129").Value = 2059.2222
$ws.Range("I129").Value = 472.16666
$ws.Range("K129").Value = 1416.49998
$ws.Range("M129").Value = 3583.50002
$ws.Range("H131").Value = 2249.875
$ws.Range("I131").Value = 1999.75
$ws.Range("K131").Value = 5999.25
$ws.Range("M131").Value = -959.25
$ws.Range("H134").Value = 2598.8125
$ws.Range("I134").Value = 2598.8125
$ws.Range("K134").Value = 7796.4375
$ws.Range("M134").Value = -2726.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 10376.556
$ws.Range("I99").Value = 5548.625
$ws.Range("J99").Value = 49000
$ws.Range("K99").Value = 5548.625
$ws.Range("L99").Value = 49000
$ws.Range("M99").Value = -3302.625
$ws.Range("N99").Value = -53492
$ws.Range("H132").Value = 31253478
$ws.Range("I132").Value = 41667970
$ws.Range("J132").Value = 10014
$ws.Range("K132").Value = 125003910
$ws.Range("L132").Value = 30042
$ws.Range("M132").Value = -125001380
$ws.Range("N132").Value = -35102
$ws.Range("H136").Value = 59146.43
$ws.Range("J136").Value = 59146.43
$ws.Range("L136").Value = 177439.29
$ws.Range("N136").Value = -182539.29

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2482.862
$ws.Range("I61").Value = 1889.037
$ws.Range("K61").Value = 1889.037
$ws.Range("M61").Value = -1687.037
$ws.Range("H68").Value = 67877.07000000001
$ws.Range("I68").Value = 1297.5
$ws.Range("J68").Value = 334195.34
$ws.Range("K68").Value = 1297.5
$ws.Range("L68").Value = 334195.34
$ws.Range("M68").Value = -548.5
$ws.Range("N68").Value = -335693.34
$ws.Range("H71").Value = 67877.07000000001
$ws.Range("I71").Value = 1297.5
$ws.Range("J71").Value = 334195.34
$ws.Range("K71").Value = 6487.5
$ws.Range("L71").Value = 1670976.7
$ws.Range("M71").Value = -2743.5
$ws.Range("N71").Value = -1678464.7
$ws.Range("H93").Value = 611.3889
$ws.Range("I93").Value = 551.2143
$ws.Range("K93").Value = 551.2143
$ws.Range("M93").Value = 696.7857
$ws.Range("H113").Value = 2482.862
$ws.Range("I113").Value = 1889.037
$ws.Range("K113").Value = 1889.037
$ws.Range("M113").Value = 280.963

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 42330
$ws.Range("J98").Value = 42330
$ws.Range("L98").Value = 42330
$ws.Range("N98").Value = -48320
$ws.Range("H132").Value = 41683920
$ws.Range("I132").Value = 50009704
$ws.Range("K132").Value = 150029112
$ws.Range("M132").Value = -150026582
$ws.Range("H139").Value = 159142.14
$ws.Range("J139").Value = 159142.14
$ws.Range("L139").Value = 159142.14
$ws.Range("N139").Value = -169422.14
